$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

$ws.Range("D2").Value = -7.759
$ws.Range("D5").Value = -7.855999999999999
$ws.Range("D6").Value = -7.797
$ws.Range("D8").Value = -7.81
$ws.Range("B11").Value = 6.4
$ws.Range("A12").Value = -21.448
$ws.Range("C14").Value = -12.697
$ws.Range("D17").Value = -8.102999999999998
$ws.Range("C19").Value = -12.385
$ws.Range("B23").Value = 7.779000000000001
$ws.Range("C24").Value = -12.687
$ws.Range("A27").Value = -21.632
$ws.Range("D27").Value = -7.908000000000001
$ws.Range("B28").Value = 5.203
$ws.Range("A32").Value = -20.982
$ws.Range("B32").Value = 8.013
$ws.Range("B34").Value = 7.049000000000001
$ws.Range("A36").Value = -20.519
$ws.Range("A38").Value = -20.528
$ws.Range("C38").Value = -11.64
$ws.Range("C41").Value = -11.991
$ws.Range("B42").Value = 7.761
$ws.Range("A46").Value = -21.651
$ws.Range("B49").Value = 6.607000000000001
$ws.Range("C52").Value = -11.632
$ws.Range("A54").Value = -20.962
$ws.Range("B54").Value = 6.215999999999999
$ws.Range("A55").Value = -22.18
$ws.Range("D55").Value = -7.848000000000001
$ws.Range("A56").Value = -21.461
$ws.Range("A67").Value = -21.359
$ws.Range("A69").Value = -21.323
$ws.Range("D70").Value = -7.007
$ws.Range("A72").Value = -21.194
$ws.Range("C72").Value = -12.648
$ws.Range("B78").Value = 7.823
$ws.Range("C78").Value = -11.868
$ws.Range("B80").Value = 7.503
$ws.Range("D80").Value = -7.583
$ws.Range("A83").Value = -21.229
$ws.Range("C83").Value = -13.056
$ws.Range("C85").Value = -11.742
$ws.Range("A86").Value = -21.429
$ws.Range("C86").Value = -13.293
$ws.Range("C90").Value = -10.589
$ws.Range("A91").Value = -21.056
$ws.Range("A93").Value = -21.338
$ws.Range("D95").Value = -7.68
$ws.Range("C96").Value = -11.854
$ws.Range("B97").Value = 5.543000000000001
$ws.Range("D98").Value = -7.802000000000001
$ws.Range("A99").Value = -20.861
$ws.Range("B99").Value = 6.287999999999999
$ws.Range("B101").Value = 5.829000000000001
$ws.Range("D102").Value = -7.965000000000001
$ws.Range("C103").Value = -12.997
$ws.Range("A104").Value = -21.437
